$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Tighten the iterative-calculation convergence threshold (iterateDelta 0.0001 -> 0.001)
$excel.Iterative = $true
$excel.MaxChange = 0.001
$excel.MaxIterations = 100
$excel.Iterative = $false

# 2. Translate the Hungarian "összesen" label to English "sum" in the office subtotal template
$ws.Range("C4").Value = "{activity.embed(`${activity.office} sum:`)}"
